$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet tracks a rolling window of trading days. "Today" rolled forward,
# so 11 newer rows (2018-03-19 .. 2018-05-09) need to be inserted right
# after the header, and the 11 oldest rows (2017-06-02 .. 2017-06-19 area)
# fall off the bottom automatically because we only ever keep a fixed
# lookback. First, shift the existing 62 data rows (A2:E63) down by 11
# rows (to A13:E74) using Copy, which relocates values/shared-strings/
# styles faithfully without introducing any new style records.
$ws.Range("A2:E63").Copy($ws.Range("A13"))

# The workbook's convention marks only the very first data row's column-A
# flag as literal text "0" (everything else is the number 0); after the
# Copy above, the row that used to be "row 2" (now row 13) still carries
# that text flavour, so normalize it back to a plain number.
$ws.Range("A13").Value = 0

# New rows of data to place into the now-empty A2:E12 gap.
$newDates = @(
    "2018-05-09",
    "2018-05-08",
    "2018-05-07",
    "2018-05-04",
    "2018-05-03",
    "2018-05-02",
    "2018-03-23",
    "2018-03-22",
    "2018-03-21",
    "2018-03-20",
    "2018-03-19"
)
$newC = @(517035.99, 517035.99, 517035.99, 517035.99, 517035.99, 517035.99, 517035.99, 517035.99, 517035.99, 517035.99, 517035.99)
$newD = @(3825.86945641, 4437.8465731, 4219.48144735, 3603.79562508, 4074.55593171, 3851.82859712, 6353.12062885, 4237.90236978, 4815.35781937, 4022.15216705, 3969.02889624)
$newE = @(0.739961923426259, 0.8583244994415186, 0.8160904712552022, 0.6970105939975281, 0.7880604078857257, 0.7449826843040462, 1.228757910808105, 0.8196532643269185, 0.9313390000897228, 0.7779249887517502, 0.7676504098370406)

for ($i = 0; $i -lt $newDates.Length; $i++) {
    $r = 2 + $i
    # Route the date text through a formula then paste-special as values so
    # it lands as a literal shared string (t="s") instead of Excel's
    # auto-date-parsing turning "2018-05-09" into a date serial number —
    # matching how the existing date column is stored.
    $ws.Range("B$r").Formula = '="' + $newDates[$i] + '"'
    $ws.Range("C$r").Value = $newC[$i]
    $ws.Range("D$r").Value = $newD[$i]
    $ws.Range("E$r").Value = $newE[$i]
}

$bRange = $ws.Range("B2:B12")
$bRange.Copy()
$bRange.PasteSpecial(-4163)
